$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their text (string) format, matching the
# original inline-string cell type, instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.293.55'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.524.37'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.49'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.95'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.522.94'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("E10").Value = '  -4.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.04'
$ws.Range("E11").Value = '  +2.31%  '
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.121.80'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000207'
$ws.Range("E14").Value = '  -3.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.40'
$ws.Range("E15").Value = '  -4.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.520.90'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.358.09'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.72'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '426.77'
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.37'
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.659.25'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000119'
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("E28").Value = '  -5.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.03'
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("E33").Value = '  -7.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.28'
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.507.88'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.62'
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.79'
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '170.60'
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.19'
$ws.Range("E43").Value = '  -4.73%  '
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.89'
$ws.Range("E45").Value = '  -9.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.46'
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("E47").Value = '  -8.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.91'
$ws.Range("E48").Value = '  -11.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.41'
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.18'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.946'
$ws.Range("E51").Value = '  -4.09%  '
